$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header block: fill in POSITION (NURSE) and UNIT (ONT) ---
$ws.Range("B3").Value = "NURSE"
$ws.Range("F4").Value = "ONT"

# --- Leave card body: continue the monthly PERIOD dates down through row 55 ---
$monthlyDates = @(44958,44986,45017,45047,45078,45108,45139,45170,45200,45231, `
                  45261,45292,45323,45352,45383,45413,45444,45474,45505,45536, `
                  45566,45597,45627,45658,45689,45717,45748,45778,45809,45839, `
                  45870,45901,45931,45962,45992)

$r = 21
foreach ($d in $monthlyDates) {
    $ws.Cells.Item($r, 1).Value = $d
    $r = $r + 1
}

# --- Row 21: SL earned this period ---
$ws.Range("C21").Value = 1.25

# --- Row 22: CASUAL->SL leave usage entry (3 days SL used 2/26-28/2023) ---
$ws.Range("B22").Value = "SL(3-0-0)"
$ws.Range("C22").Value = 1.25
$ws.Range("H22").Value = 3
$ws.Range("K22").Value = "2/26-28/2023"

# --- Restore the selection that was left active on save ---
$ws.Activate()
$ws.Range("B23").Select()
